$wb = $excel.ActiveWorkbook

# Update the MarketCap value on the AddOpportunity sheet (row 2, col AA)
$ws1 = $wb.Worksheets.Item("AddOpportunity")
$ws1.Range("AA2").Value = "10000.0"

# Make AddOpportunity the active/selected sheet (was ValuationPeriod)
$ws1.Activate()

# Update the active selection on AddOpportunity to AA3
$null = $ws1.Range("AA3").Select()
